$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for "Reset Passwod functionality" right before the
# --- blank separator row (old row 28), shifting the "Tehnical" section down by one.
$ws.Rows.Item(28).Insert()
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "Reset Passwod functionality"

# --- Append a new last row (42) for "Make swagger work with jwtToken header"
# --- Copy the formatting of the row above (41) so the new row matches the
# --- existing table styling (font/fill/border of columns B and C).
$ws.Range("A41:C41").Copy()
$ws.Range("A42:C42").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A42").Value = 11
$ws.Range("B42").Value = "Make swagger work with jwtToken header"

# --- Restore/update the view selection to match the final state
$ws.Range("D42").Select()
